# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The underlying worker database was refreshed, which re-sorted the workers
# listed in the account-statement table (Hoja1, rows 16-22). Each worker's
# document number, name and the related "Periodo Mora" / "Valor Mora"
# figures (columns C, D, F and G) now appear in a new row order; columns B
# (document type) and E (periodo) are identical for every worker and stay
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New contents for rows 16-22: DocNumber (C), Name (D), PeriodoMora (F),
# ValorMora (G).
$updates = @(
    @{ Row = 16; C = "1047380341"; D = "MERY ISABEL ARZUZA TORRES";     F = 12887; G = 644350 },
    @{ Row = 17; C = "1128045187"; D = "CAMILO ANDRES PEREZ LUJAN";     F = 16000; G = 800000 },
    @{ Row = 18; C = "1082878492"; D = "JOSE DAVID MONTES FERRADANEZ"; F = 12887; G = 644350 },
    @{ Row = 19; C = "73074492";   D = "JULIO CESAR MEZA URZOLA";       F = 12320; G = 616000 },
    @{ Row = 20; C = "1143334730"; D = "JORGE LEONARDO VALDEZ TAPIA";   F = 12320; G = 616000 },
    @{ Row = 21; C = "73137385";   D = "EDINSON MANUEL VENECIA PITALUA"; F = 13789; G = 689455 },
    @{ Row = 22; C = "1047477071"; D = "EDDY MANUEL GUTIERREZ OSORIO";  F = 12320; G = 616000 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C   # Column C - N Doc Trabajador
    $ws.Cells.Item($u.Row, 4).Value = $u.D   # Column D - Nombre Trabajador
    $ws.Cells.Item($u.Row, 6).Value = $u.F   # Column F - Periodo Mora (count)
    $ws.Cells.Item($u.Row, 7).Value = $u.G   # Column G - Valor Mora
}
